$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) RF1 paragraph: rewrite the sentence text.
# ---------------------------------------------------------------------------
$oldRF1 = "RF1. Registrar un restaurante con un nombre, NIT y el nombre del administrador. "
$newRF1 = "RF1. Registrar un restaurante. Para esto se necesita el nombre del restaurante, el NIT y el nombre del administrador. Ya sea que se registre exitosamente el restaurante o que ocurra algún error se le informa al usuario lo que ocurrió. "

$found = $d.Content.Find.Execute($oldRF1, $true, $false, $false, $false, $false, $true, 1, $false, $newRF1, 2)

# Re-split the (now single) run for this paragraph into the same run
# boundaries the human edit produced, by toggling a formatting property on
# and back off for each slice - this creates separate <w:r> elements while
# leaving every run's effective formatting identical (lang=es-MX only).
$rf1Para = $d.Paragraphs(2)
$rf1Start = $rf1Para.Range.Start

$boundaries = @(29, 52, 53, 55, 62, 78, 79, 82, 118, 235)
$prev = 0
foreach ($b in $boundaries) {
    $chunk = $d.Range($rf1Start + $prev, $rf1Start + $b)
    $chunk.Bold = $true
    $chunk2 = $d.Range($rf1Start + $prev, $rf1Start + $b)
    $chunk2.Bold = $false
    $prev = $b
}

# ---------------------------------------------------------------------------
# 2) RF20 / RF21 paragraphs: merge the runs that make up the "RFxx. " prefix
#    and the "con informacion de ..." suffix (leaving the cvs proofed word
#    runs untouched), matching the up-to-date canonical run layout.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("RF20. Importar datos de un archivo ", $true, $false, $false, $false, $false, $true, 1, $false, "RF20. Importar datos de un archivo ", 2)
$null = $d.Content.Find.Execute(" con información de productos", $true, $false, $false, $false, $false, $true, 1, $false, " con información de productos", 2)

$null = $d.Content.Find.Execute("RF21. Importar datos de un archivo ", $true, $false, $false, $false, $false, $true, 1, $false, "RF21. Importar datos de un archivo ", 2)
$null = $d.Content.Find.Execute(" con información de pedidos", $true, $false, $false, $false, $false, $true, 1, $false, " con información de pedidos", 2)
